# Recall-results workbook update:
# - add RNN / bi-LSTM style recall column (G) with bold, centered, 2-decimal values
# - add delta rows comparing E/F/G against D, and G against F
# - add a second sheet "politics" with a domain/recall@1 header + hyperlink

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Header row: center-align C2:G2 (keeps existing bold/border formatting) ----
$ws1.Range("C2:G2").HorizontalAlignment = -4108   # xlCenter

# ---- Existing numeric block C3:F5: switch to 2-decimal format, centered ----
$ws1.Range("C3:F5").NumberFormat = "0.00"
$ws1.Range("C3:F5").HorizontalAlignment = -4108   # xlCenter

# ---- New column G: Bi-LSTM-with-RD recall values, bold + centered + 2 decimals ----
$ws1.Range("G3").Value = 0.62343752399999997
$ws1.Range("G4").Value = 0.70390623799999996
$ws1.Range("G5").Value = 0.83984375
$ws1.Range("G3:G5").NumberFormat = "0.00"
$ws1.Range("G3:G5").HorizontalAlignment = -4108   # xlCenter
$ws1.Range("G3:G5").Font.Bold = $true

# ---- Delta rows 7:9 -> recall(E/F/G) - recall(D) ----
$ws1.Range("E7").Formula = '=E3-$D3'
$ws1.Range("F7:G7").Formula = '=F3-$D3'
$ws1.Range("E8").Formula = '=E4-$D4'
$ws1.Range("F8:G8").Formula = '=F4-$D4'
$ws1.Range("E9").Formula = '=E5-$D5'
$ws1.Range("F9:G9").Formula = '=F5-$D5'
$ws1.Range("E7:G9").NumberFormat = "0.00"

# ---- Delta rows 11:13 -> recall(G) - recall(F) ----
$ws1.Range("G11").Formula = '=G3-F3'
$ws1.Range("G12:G13").Formula = '=G4-F4'
$ws1.Range("G11:G13").NumberFormat = "0.00"

# ---- Selection / view bookkeeping on Sheet1 (matches the saved file) ----
$ws1.Range("G11:G13").Select()

# ---- Add the new "politics" sheet after Sheet1 ----
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "politics"

$ws2.Range("B2").Value = "domain"
$ws2.Range("C2").Value = "recall@1"
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://d.docs.live.net/af91f69c463c54d7/research/twconv/twconvrsu/twconvrsu/results/recall_results.xlsx") | Out-Null

$ws2.Range("C3").Select()
